$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1264
$ws1.Range("F5").Value = 1097
$ws1.Range("F6").Value = 14093
$ws1.Range("F7").Value = 15629
$ws1.Range("F11").Value = 189
$ws1.Range("F19").Value = 28
$ws1.Range("F23").Value = 6115
$ws1.Range("F26").Value = 5548
$ws1.Range("F27").Value = 74
$ws1.Range("F28").Value = 141
$ws1.Range("F29").Value = 117
$ws1.Range("F30").Value = 4523

# Sheet "全部类型" (sheet4.xml)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1264
$ws4.Range("F5").Value = 1097
$ws4.Range("F6").Value = 14093
$ws4.Range("F7").Value = 15629
$ws4.Range("F11").Value = 189
$ws4.Range("F19").Value = 28
$ws4.Range("F24").Value = 6115
$ws4.Range("F27").Value = 5548
$ws4.Range("F28").Value = 74
$ws4.Range("F29").Value = 141
$ws4.Range("F30").Value = 117
$ws4.Range("F31").Value = 4523
